$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040379440627602
$ws.Cells.Item(2, 4).Value = 1.042043250366263
$ws.Cells.Item(2, 5).Value = 1.038533635975398
$ws.Cells.Item(2, 6).Value = 1.039131590117458
$ws.Cells.Item(2, 9).Value = 1.041206215823784
$ws.Cells.Item(2, 10).Value = 1.045466364141586
$ws.Cells.Item(2, 11).Value = 1.04482105367315
$ws.Cells.Item(2, 12).Value = 1.04132139708542
$ws.Cells.Item(2, 13).Value = 1.041917649495794
$ws.Cells.Item(2, 14).Value = 1.046951046295335
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042162155873063
$ws.Cells.Item(3, 4).Value = 1.042912896762855
$ws.Cells.Item(3, 5).Value = 1.040097250559813
$ws.Cells.Item(3, 6).Value = 1.041551369601292
$ws.Cells.Item(3, 9).Value = 1.041723130406639
$ws.Cells.Item(3, 10).Value = 1.046890516222773
$ws.Cells.Item(3, 11).Value = 1.045501296729258
$ws.Cells.Item(3, 12).Value = 1.042693055424887
$ws.Cells.Item(3, 13).Value = 1.044143344862835
$ws.Cells.Item(3, 14).Value = 1.04837722083583
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043310810598504
$ws.Cells.Item(4, 4).Value = 1.043473443725672
$ws.Cells.Item(4, 5).Value = 1.041104226440361
$ws.Cells.Item(4, 6).Value = 1.043111546096081
$ws.Cells.Item(4, 9).Value = 1.042054381371893
$ws.Cells.Item(4, 10).Value = 1.047806923689169
$ws.Cells.Item(4, 11).Value = 1.045938654614618
$ws.Cells.Item(4, 12).Value = 1.043575351235984
$ws.Cells.Item(4, 13).Value = 1.045577658461863
$ws.Cells.Item(4, 14).Value = 1.049294929705942
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.043792560090605
$ws.Cells.Item(5, 4).Value = 1.043708585062884
$ws.Cells.Item(5, 5).Value = 1.041526433671858
$ws.Cells.Item(5, 6).Value = 1.043766144386214
$ws.Cells.Item(5, 9).Value = 1.042192874002964
$ws.Cells.Item(5, 10).Value = 1.048190976936584
$ws.Cells.Item(5, 11).Value = 1.046121855806687
$ws.Cells.Item(5, 12).Value = 1.043945028877865
$ws.Cells.Item(5, 13).Value = 1.046179275008825
$ws.Cells.Item(5, 14).Value = 1.049679528353019
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.043873381490054
$ws.Cells.Item(6, 4).Value = 1.04374803648343
$ws.Cells.Item(6, 5).Value = 1.041597258723437
$ws.Cells.Item(6, 6).Value = 1.04387597914613
$ws.Cells.Item(6, 9).Value = 1.042216082860683
$ws.Cells.Item(6, 10).Value = 1.048255391109934
$ws.Cells.Item(6, 11).Value = 1.046152577365547
$ws.Cells.Item(6, 12).Value = 1.044007027255374
$ws.Cells.Item(6, 13).Value = 1.046280209782325
$ws.Cells.Item(6, 14).Value = 1.049744034001884
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.043317252225802
$ws.Cells.Item(7, 4).Value = 1.04347658769969
$ws.Cells.Item(7, 5).Value = 1.041109872388305
$ws.Cells.Item(7, 6).Value = 1.043120297929094
$ws.Cells.Item(7, 9).Value = 1.042056234914229
$ws.Cells.Item(7, 10).Value = 1.047812060133277
$ws.Cells.Item(7, 11).Value = 1.045941105157104
$ws.Cells.Item(7, 12).Value = 1.043580295730899
$ws.Cells.Item(7, 13).Value = 1.045585702624684
$ws.Cells.Item(7, 14).Value = 1.049300073444389
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040982942595784
$ws.Cells.Item(8, 4).Value = 1.042337605782759
$ws.Cells.Item(8, 5).Value = 1.039063070605475
$ws.Cells.Item(8, 6).Value = 1.039950544007233
$ws.Cells.Item(8, 9).Value = 1.041381582227749
$ws.Cells.Item(8, 10).Value = 1.04594873397511
$ws.Cells.Item(8, 11).Value = 1.045051530267052
$ws.Cells.Item(8, 12).Value = 1.041786056113719
$ws.Cells.Item(8, 13).Value = 1.04267106766305
$ws.Cells.Item(8, 14).Value = 1.047434101149352
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.036831136470355
$ws.Cells.Item(9, 4).Value = 1.040313596164816
$ws.Cells.Item(9, 5).Value = 1.035418743571324
$ws.Cells.Item(9, 6).Value = 1.034320619654196
$ws.Cells.Item(9, 9).Value = 1.040167706572623
$ws.Cells.Item(9, 10).Value = 1.042625245837122
$ws.Cells.Item(9, 11).Value = 1.043462166349285
$ws.Cells.Item(9, 12).Value = 1.03858321864882
$ws.Cells.Item(9, 13).Value = 1.03748868553639
$ws.Cells.Item(9, 14).Value = 1.044105893276998
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034035902745237
$ws.Cells.Item(10, 4).Value = 1.038952358943317
$ws.Cells.Item(10, 5).Value = 1.032962592037631
$ws.Cells.Item(10, 6).Value = 1.030535141297969
$ws.Cells.Item(10, 9).Value = 1.039341140584799
$ws.Cells.Item(10, 10).Value = 1.040381375511767
$ws.Cells.Item(10, 11).Value = 1.042387442145294
$ws.Cells.Item(10, 12).Value = 1.03641907813085
$ws.Cells.Item(10, 13).Value = 1.034000379474159
$ws.Cells.Item(10, 14).Value = 1.041858836398406
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032818706235944
$ws.Cells.Item(11, 4).Value = 1.03835999680404
$ws.Cells.Item(11, 5).Value = 1.031892445227474
$ws.Cells.Item(11, 6).Value = 1.028887826231901
$ws.Cells.Item(11, 9).Value = 1.038979014626786
$ws.Cells.Item(11, 10).Value = 1.039402777326903
$ws.Cells.Item(11, 11).Value = 1.041918373450965
$ws.Cells.Item(11, 12).Value = 1.035474843265292
$ws.Cells.Item(11, 13).Value = 1.032481497610621
$ws.Cells.Item(11, 14).Value = 1.040878848491872
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032365527868643
$ws.Cells.Item(12, 4).Value = 1.038139516216319
$ws.Cells.Item(12, 5).Value = 1.031493925124522
$ws.Cells.Item(12, 6).Value = 1.028274665902957
$ws.Cells.Item(12, 9).Value = 1.038843861969545
$ws.Cells.Item(12, 10).Value = 1.039038207706704
$ws.Cells.Item(12, 11).Value = 1.041743574183637
$ws.Cells.Item(12, 12).Value = 1.035123014202711
$ws.Cells.Item(12, 13).Value = 1.03191600951936
$ws.Cells.Item(12, 14).Value = 1.040513761140995
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032462784494186
$ws.Cells.Item(13, 4).Value = 1.038186830611666
$ws.Cells.Item(13, 5).Value = 1.031579455612091
$ws.Cells.Item(13, 6).Value = 1.028406249272662
$ws.Cells.Item(13, 9).Value = 1.038872881919837
$ws.Cells.Item(13, 10).Value = 1.039116458163361
$ws.Cells.Item(13, 11).Value = 1.041781094998776
$ws.Cells.Item(13, 12).Value = 1.035198532821179
$ws.Cells.Item(13, 13).Value = 1.032037368500472
$ws.Cells.Item(13, 14).Value = 1.040592122722277
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032781268113414
$ws.Cells.Item(14, 4).Value = 1.038341781064582
$ws.Cells.Item(14, 5).Value = 1.031859524388496
$ws.Cells.Item(14, 6).Value = 1.02883716846661
$ws.Cells.Item(14, 9).Value = 1.038967856036604
$ws.Cells.Item(14, 10).Value = 1.03937266395704
$ws.Cells.Item(14, 11).Value = 1.041903936096198
$ws.Cells.Item(14, 12).Value = 1.03544578351017
$ws.Cells.Item(14, 13).Value = 1.032434781113868
$ws.Cells.Item(14, 14).Value = 1.040848692357569
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03297735537809
$ws.Cells.Item(15, 4).Value = 1.038437191111821
$ws.Cells.Item(15, 5).Value = 1.032031948075304
$ws.Cells.Item(15, 6).Value = 1.029102501789412
$ws.Cells.Item(15, 9).Value = 1.039026287244967
$ws.Cells.Item(15, 10).Value = 1.039530377644849
$ws.Cells.Item(15, 11).Value = 1.041979547227781
$ws.Cells.Item(15, 12).Value = 1.035597976564295
$ws.Cells.Item(15, 13).Value = 1.032679465572777
$ws.Cells.Item(15, 14).Value = 1.041006630016905
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03411653741385
$ws.Cells.Item(16, 4).Value = 1.038991609364149
$ws.Cells.Item(16, 5).Value = 1.033033472551997
$ws.Cells.Item(16, 6).Value = 1.030644292005632
$ws.Cells.Item(16, 9).Value = 1.039365084042868
$ws.Cells.Item(16, 10).Value = 1.040446172516615
$ws.Cells.Item(16, 11).Value = 1.042418493802808
$ws.Cells.Item(16, 12).Value = 1.03648159122738
$ws.Cells.Item(16, 13).Value = 1.034101001611541
$ws.Cells.Item(16, 14).Value = 1.041923725422433
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034829263791715
$ws.Cells.Item(17, 4).Value = 1.039338588132637
$ws.Cells.Item(17, 5).Value = 1.033659912758833
$ws.Cells.Item(17, 6).Value = 1.03160919491071
$ws.Cells.Item(17, 9).Value = 1.039576466661673
$ws.Cells.Item(17, 10).Value = 1.041018738248279
$ws.Cells.Item(17, 11).Value = 1.042692834546846
$ws.Cells.Item(17, 12).Value = 1.037033928661702
$ws.Cells.Item(17, 13).Value = 1.034990409025478
$ws.Cells.Item(17, 14).Value = 1.04249710426311
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.035244326904399
$ws.Cells.Item(18, 4).Value = 1.03954069221572
$ws.Cells.Item(18, 5).Value = 1.034024668053535
$ws.Cells.Item(18, 6).Value = 1.03217122050872
$ws.Cells.Item(18, 9).Value = 1.039699356284876
$ws.Cells.Item(18, 10).Value = 1.041352033659288
$ws.Cells.Item(18, 11).Value = 1.04285249593366
$ws.Cells.Item(18, 12).Value = 1.037355409752084
$ws.Cells.Item(18, 13).Value = 1.03550837558872
$ws.Cells.Item(18, 14).Value = 1.042830872991832
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.035385741954941
$ws.Cells.Item(19, 4).Value = 1.039609556811292
$ws.Cells.Item(19, 5).Value = 1.034148932965505
$ws.Cells.Item(19, 6).Value = 1.032362724770513
$ws.Cells.Item(19, 9).Value = 1.039741189837969
$ws.Cells.Item(19, 10).Value = 1.041465565461628
$ws.Cells.Item(19, 11).Value = 1.042906876083546
$ws.Cells.Item(19, 12).Value = 1.037464910531111
$ws.Cells.Item(19, 13).Value = 1.035684852728422
$ws.Cells.Item(19, 14).Value = 1.042944566022353
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034752863281512
$ws.Cells.Item(20, 4).Value = 1.039301389911004
$ws.Cells.Item(20, 5).Value = 1.033592767693539
$ws.Cells.Item(20, 6).Value = 1.031505751581566
$ws.Cells.Item(20, 9).Value = 1.039553829397945
$ws.Cells.Item(20, 10).Value = 1.040957377102073
$ws.Cells.Item(20, 11).Value = 1.042663437368242
$ws.Cells.Item(20, 12).Value = 1.036974739419278
$ws.Cells.Item(20, 13).Value = 1.034895068086417
$ws.Cells.Item(20, 14).Value = 1.042435655977038
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032687512053766
$ws.Cells.Item(21, 4).Value = 1.038296164543663
$ws.Cells.Item(21, 5).Value = 1.031777079478734
$ws.Cells.Item(21, 6).Value = 1.028710309047026
$ws.Cells.Item(21, 9).Value = 1.038939906350789
$ws.Cells.Item(21, 10).Value = 1.039297247604974
$ws.Cells.Item(21, 11).Value = 1.041867778163546
$ws.Cells.Item(21, 12).Value = 1.03537300485762
$ws.Cells.Item(21, 13).Value = 1.032317789448394
$ws.Cells.Item(21, 14).Value = 1.040773168905633
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031382811084822
$ws.Cells.Item(22, 4).Value = 1.037661525718347
$ws.Cells.Item(22, 5).Value = 1.030629570145508
$ws.Cells.Item(22, 6).Value = 1.026945308321942
$ws.Cells.Item(22, 9).Value = 1.038550183894283
$ws.Cells.Item(22, 10).Value = 1.038247227887235
$ws.Cells.Item(22, 11).Value = 1.041364234559528
$ws.Cells.Item(22, 12).Value = 1.034359564626221
$ws.Cells.Item(22, 13).Value = 1.030689764380935
$ws.Cells.Item(22, 14).Value = 1.039721658039457
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03207504903004
$ws.Cells.Item(23, 4).Value = 1.037998210771936
$ws.Cells.Item(23, 5).Value = 1.031238455867466
$ws.Cells.Item(23, 6).Value = 1.027881684968844
$ws.Cells.Item(23, 9).Value = 1.038757139356862
$ws.Cells.Item(23, 10).Value = 1.03880446219565
$ws.Cells.Item(23, 11).Value = 1.041631486771321
$ws.Cells.Item(23, 12).Value = 1.034897420130655
$ws.Cells.Item(23, 13).Value = 1.031553544837553
$ws.Cells.Item(23, 14).Value = 1.040279683684505
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034787387416011
$ws.Cells.Item(24, 4).Value = 1.039318199059763
$ws.Cells.Item(24, 5).Value = 1.033623109624763
$ws.Cells.Item(24, 6).Value = 1.03155249559602
$ws.Cells.Item(24, 9).Value = 1.039564059458129
$ws.Cells.Item(24, 10).Value = 1.040985105639894
$ws.Cells.Item(24, 11).Value = 1.042676721790543
$ws.Cells.Item(24, 12).Value = 1.037001486615172
$ws.Cells.Item(24, 13).Value = 1.034938151053211
$ws.Cells.Item(24, 14).Value = 1.042463423892563
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.037909188752957
$ws.Cells.Item(25, 4).Value = 1.040838912248529
$ws.Cells.Item(25, 5).Value = 1.036365479931365
$ws.Cells.Item(25, 6).Value = 1.035781591283816
$ws.Cells.Item(25, 9).Value = 1.040484539377347
$ws.Cells.Item(25, 10).Value = 1.043489325231669
$ws.Cells.Item(25, 11).Value = 1.043875689414276
$ws.Cells.Item(25, 12).Value = 1.039416233132108
$ws.Cells.Item(25, 13).Value = 1.038834178421909
$ws.Cells.Item(25, 14).Value = 1.044971199763397
